# Weekly fruit/vegetable price update: insert 2 new rows of data at the top
# of the "Limón" price-history block (rows 346-356), pushing the existing
# rows down to 348-358, then populate the two newly inserted rows
# (346 and 347) with this week's observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 346 (existing row 346 and everything
# below shifts down by two, so old 346..356 become 348..358).
$ws.Rows.Item(346).Resize(2).Insert()

# ---- New row 346: Sutil De Gase / Primera, $/caja 18 kilos ----
$ws.Cells.Item(346, 1).Value  = 1
$ws.Cells.Item(346, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(346, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(346, 4).Value  = 44939
$ws.Cells.Item(346, 5).Value  = 15
$ws.Cells.Item(346, 6).Value  = "Fruta"
$ws.Cells.Item(346, 7).Value  = 100102
$ws.Cells.Item(346, 8).Value  = "Cítricos"
$ws.Cells.Item(346, 9).Value  = 100102003
$ws.Cells.Item(346, 10).Value = "Limón"
$ws.Cells.Item(346, 11).Value = "Sutil De Gase"
$ws.Cells.Item(346, 12).Value = "Primera"
$ws.Cells.Item(346, 13).Value = 250
$ws.Cells.Item(346, 14).Value = 28000
$ws.Cells.Item(346, 15).Value = 29000
$ws.Cells.Item(346, 16).Value = 28600
$ws.Cells.Item(346, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(346, 18).Value = "Perú"
$ws.Cells.Item(346, 19).Value = 1589
$ws.Cells.Item(346, 20).Value = 18

# ---- New row 347: Tahití / Primera, $/caja 24 kilos ----
$ws.Cells.Item(347, 1).Value  = 1
$ws.Cells.Item(347, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(347, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(347, 4).Value  = 44939
$ws.Cells.Item(347, 5).Value  = 15
$ws.Cells.Item(347, 6).Value  = "Fruta"
$ws.Cells.Item(347, 7).Value  = 100102
$ws.Cells.Item(347, 8).Value  = "Cítricos"
$ws.Cells.Item(347, 9).Value  = 100102003
$ws.Cells.Item(347, 10).Value = "Limón"
$ws.Cells.Item(347, 11).Value = "Tahití"
$ws.Cells.Item(347, 12).Value = "Primera"
$ws.Cells.Item(347, 13).Value = 350
$ws.Cells.Item(347, 14).Value = 30000
$ws.Cells.Item(347, 15).Value = 31000
$ws.Cells.Item(347, 16).Value = 30429
$ws.Cells.Item(347, 17).Value = "$/caja 24 kilos"
$ws.Cells.Item(347, 18).Value = "Perú"
$ws.Cells.Item(347, 19).Value = 1268
$ws.Cells.Item(347, 20).Value = 24
